# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Tue Apr 30 12:42:12 UTC 2024 with GitHub Actions"
#
# Most rows just get refreshed Price / Volume(1h) figures. Rows 24/25 and
# 40/41 swapped their ranking order, so Coin / Link / Price / Volume(1h) are
# all rewritten for those four rows.
#
# Price text that looks like a plain decimal (e.g. "569.30") is written with a
# leading apostrophe so Excel keeps it as text (matching the source data, which
# stores prices/volumes as formatted strings, not numbers) instead of silently
# parsing it into a float and dropping the trailing zero. The Style reset right
# after puts the cell back on the workbook's default "Normal" style so the
# quote-prefix marker Excel adds for that doesn't stick around as a formatting
# change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.240.50'
$ws.Range("E2").Value = '  -1.95%  '

# Row 3
$ws.Range("D3").Value = '3.009.94'
$ws.Range("E3").Value = '  -4.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").Value = '''569.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.31%  '

# Row 6
$ws.Range("D6").Value = '''128.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.71%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '3.010.08'
$ws.Range("E8").Value = '  -4.36%  '

# Row 9
$ws.Range("D9").Value = '''0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.23%  '

# Row 10
$ws.Range("D10").Value = '''0.135'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.71%  '

# Row 11
$ws.Range("D11").Value = '''5.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.22%  '

# Row 12
$ws.Range("E12").Value = '  -5.21%  '

# Row 13
$ws.Range("E13").Value = '  -3.66%  '

# Row 14
$ws.Range("D14").Value = '''32.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.42%  '

# Row 15
$ws.Range("E15").Value = '  -0.43%  '

# Row 16
$ws.Range("D16").Value = '3.503.51'
$ws.Range("E16").Value = '  -4.44%  '

# Row 17
$ws.Range("D17").Value = '61.270.00'
$ws.Range("E17").Value = '  -1.72%  '

# Row 18
$ws.Range("D18").Value = '3.010.99'
$ws.Range("E18").Value = '  -4.83%  '

# Row 19
$ws.Range("D19").Value = '''6.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.57%  '

# Row 20
$ws.Range("D20").Value = '''439.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.94%  '

# Row 21
$ws.Range("D21").Value = '''13.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.15%  '

# Row 22
$ws.Range("E22").Value = '  -5.43%  '

# Row 23
$ws.Range("D23").Value = '''7.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.28%  '

# Row 24
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''12.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.53%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''78.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.32%  '

# Row 26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("E27").Value = '  +0.27%  '

# Row 28
$ws.Range("D28").Value = '''2.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.86%  '

# Row 29
$ws.Range("D29").Value = '''7.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.06%  '

# Row 30
$ws.Range("D30").Value = '''1.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.39%  '

# Row 31
$ws.Range("D31").Value = '''6.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.26%  '

# Row 32
$ws.Range("D32").Value = '''25.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.27%  '

# Row 33
$ws.Range("D33").Value = '''0.0946'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.98%  '

# Row 34
$ws.Range("D34").Value = '''2.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.04%  '

# Row 35
$ws.Range("D35").Value = '''0.959'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.14%  '

# Row 36
$ws.Range("D36").Value = '''5.57'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.81%  '

# Row 37
$ws.Range("D37").Value = '''50.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.94%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0690'
$ws.Range("E38").Value = '  -0.96%  '

# Row 39
$ws.Range("E39").Value = '  -4.46%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.109'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.93%  '

# Row 41
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '''7.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.44%  '

# Row 42
$ws.Range("D42").Value = '''374.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.03%  '

# Row 43
$ws.Range("E43").Value = '  -9.87%  '

# Row 44
$ws.Range("D44").Value = '2.645.54'
$ws.Range("E44").Value = '  -5.40%  '

# Row 45
$ws.Range("E45").Value = '  -0.04%  '

# Row 46
$ws.Range("E46").Value = '  -5.23%  '

# Row 47
$ws.Range("D47").Value = '''120.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.66%  '

# Row 48
$ws.Range("D48").Value = '''33.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.43%  '

# Row 49
$ws.Range("D49").Value = '''1.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.67%  '

# Row 50
$ws.Range("E50").Value = '  -3.98%  '

# Row 51
$ws.Range("D51").Value = '''23.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.73%  '
